$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update build number
$ws.Range("C2").Value = 7977

# Reorder the changelog entries in A6:A17 (alphabetical re-sort after adding
# two new fixes: "Changed Nylium textures" and "Fixed broken mycelium texture")
$newValues = @(
    "Changed Nylium textures",
    "Changed Vex CEM",
    "Deprecated Sodium support",
    "Fixed broken mycelium texture",
    "Fixed cloud shaders",
    "Fixed mipmaps not working",
    "Fixed Rabbit CEM",
    "Improved fog shaders for 1.21.2+",
    "Refactor bat model and update texture for improved visuals",
    "Reimproved lighting (Vanilla, OptiFine)",
    "Removed OptiFine Lightmaps",
    "Updated pack format to 63 (Classic Reimagined 10 SE C2 only)"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 6 + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# Update the active selection to match the saved view
$ws.Range("I16").Select()
